$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "Uncut_Sheet_1"

# Re-point the print area so the Print_Area defined name follows the new sheet name
$ws.PageSetup.PrintArea = "A1:G42"

# Update the active selection to the next merged range (B16:C16)
$ws.Range("B16:C16").Select()
